# Amlak private - New Columns - excel import
# Rename the "کد سجام" (SEJAM code) header to "کد جام" (Jam code).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "کد جام"

# Move the active selection to I2, matching the saved view state.
$ws.Range("I2").Select()
